# Corrects mis-scaled IFRS financial figures for rows 2-9 (company_list sheet)
# of the 롯데지주 (Lotte Corporation) workbook - "error solve ifrs list".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 22248
$ws.Range("E2").Value = 1148
$ws.Range("F2").Value = 1148
$ws.Range("G2").Value = 371
$ws.Range("H2").Value = 17
$ws.Range("I2").Value = -29
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 39979
$ws.Range("L2").Value = 13777
$ws.Range("M2").Value = 26202
$ws.Range("N2").Value = 24829
$ws.Range("O2").Value = 1374
$ws.Range("P2").Value = 71
$ws.Range("Q2").Value = 941
$ws.Range("R2").Value = -1112
$ws.Range("S2").Value = -570
$ws.Range("T2").Value = 1737
$ws.Range("U2").Value = -796
$ws.Range("V2").Value = 5370
$ws.Range("W2").Value = 5.16
$ws.Range("X2").Value = 0.08
$ws.Range("Y2").Value = -0.11
$ws.Range("Z2").Value = 0.04
$ws.Range("AA2").Value = 52.58
$ws.Range("AB2").Value = 23440.8
$ws.Range("AC2").Value = -80
$ws.Range("AD2").Value = -629.51
$ws.Range("AE2").Value = 70566
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 208
$ws.Range("AH2").Value = 0.41
$ws.Range("AI2").Value = -256.71
$ws.Range("AJ2").Value = 35535000

# Row 3
$ws.Range("D3").Value = 22579
$ws.Range("E3").Value = 1445
$ws.Range("F3").Value = 1445
$ws.Range("G3").Value = 1130
$ws.Range("H3").Value = 792
$ws.Range("I3").Value = 752
$ws.Range("J3").Value = 39
$ws.Range("K3").Value = 43259
$ws.Range("L3").Value = 16233
$ws.Range("M3").Value = 27025
$ws.Range("N3").Value = 25872
$ws.Range("O3").Value = 1153
$ws.Range("P3").Value = 71
$ws.Range("Q3").Value = 2183
$ws.Range("R3").Value = -2573
$ws.Range("S3").Value = 1761
$ws.Range("T3").Value = 1532
$ws.Range("U3").Value = 651
$ws.Range("V3").Value = 7386
$ws.Range("W3").Value = 6.4
$ws.Range("X3").Value = 3.51
$ws.Range("Y3").Value = 2.97
$ws.Range("Z3").Value = 1.9
$ws.Range("AA3").Value = 60.07
$ws.Range("AB3").Value = 24718.38
$ws.Range("AC3").Value = 2117
$ws.Range("AD3").Value = 30.54
$ws.Range("AE3").Value = 72808
$ws.Range("AF3").Value = 0.89
$ws.Range("AG3").Value = 451
$ws.Range("AH3").Value = 0.7
$ws.Range("AI3").Value = 21.29
$ws.Range("AJ3").Value = 35535000

# Row 4
$ws.Range("D4").Value = 22483
$ws.Range("E4").Value = 1278
$ws.Range("F4").Value = 1278
$ws.Range("G4").Value = 962
$ws.Range("H4").Value = 752
$ws.Range("I4").Value = 698
$ws.Range("J4").Value = 54
$ws.Range("K4").Value = 39918
$ws.Range("L4").Value = 14298
$ws.Range("M4").Value = 25621
$ws.Range("N4").Value = 24426
$ws.Range("O4").Value = 1195
$ws.Range("P4").Value = 71
$ws.Range("Q4").Value = 1136
$ws.Range("R4").Value = -591
$ws.Range("S4").Value = -1360
$ws.Range("T4").Value = 1014
$ws.Range("U4").Value = 122
$ws.Range("V4").Value = 6305
$ws.Range("W4").Value = 5.68
$ws.Range("X4").Value = 3.35
$ws.Range("Y4").Value = 2.78
$ws.Range("Z4").Value = 1.81
$ws.Range("AA4").Value = 55.81
$ws.Range("AB4").Value = 25434.12
$ws.Range("AC4").Value = 1965
$ws.Range("AD4").Value = 25.78
$ws.Range("AE4").Value = 68737
$ws.Range("AF4").Value = 0.74
$ws.Range("AG4").Value = 452
$ws.Range("AH4").Value = 0.89
$ws.Range("AI4").Value = 23.01
$ws.Range("AJ4").Value = 35535000

# Row 5
$ws.Range("D5").Value = 18690
$ws.Range("E5").Value = -68
$ws.Range("F5").Value = -68
$ws.Range("G5").Value = 5174
$ws.Range("H5").Value = 4953
$ws.Range("I5").Value = 4688
$ws.Range("J5").Value = 265
$ws.Range("K5").Value = 218971
$ws.Range("L5").Value = 144072
$ws.Range("M5").Value = 74898
$ws.Range("N5").Value = 61148
$ws.Range("O5").Value = 13751
$ws.Range("P5").Value = 149
$ws.Range("Q5").Value = -6840
$ws.Range("R5").Value = 13058
$ws.Range("S5").Value = -3272
$ws.Range("T5").Value = 612
$ws.Range("U5").Value = -7452
$ws.Range("V5").Value = 98120
$ws.Range("W5").Value = -0.36
$ws.Range("X5").Value = 26.5
$ws.Range("Y5").Value = 10.96
$ws.Range("Z5").Value = 3.83
$ws.Range("AA5").Value = 192.36
$ws.Range("AB5").Value = 50440.01
$ws.Range("AC5").Value = 10602
$ws.Range("AD5").Value = 6.14
$ws.Range("AE5").Value = 100791
$ws.Range("AF5").Value = 0.65
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 73683000

# Row 6
$ws.Range("D6").Value = 72712
$ws.Range("E6").Value = 984
$ws.Range("F6").Value = 984
$ws.Range("G6").Value = 3386
$ws.Range("H6").Value = 1107
$ws.Range("I6").Value = 1851
$ws.Range("K6").Value = 276183
$ws.Range("L6").Value = 191849
$ws.Range("M6").Value = 84335
$ws.Range("N6").Value = 71130
$ws.Range("P6").Value = 212
$ws.Range("Q6").Value = -8267
$ws.Range("R6").Value = -22480
$ws.Range("S6").Value = 34136
$ws.Range("T6").Value = 2397
$ws.Range("U6").Value = -10663
$ws.Range("V6").Value = 50379
$ws.Range("W6").Value = 1.35
$ws.Range("X6").Value = 1.52
$ws.Range("Y6").Value = 2.8
$ws.Range("Z6").Value = 0.45
$ws.Range("AA6").Value = 227.48
$ws.Range("AB6").Value = 48528.98
$ws.Range("AC6").Value = 1746
$ws.Range("AD6").Value = 30.18
$ws.Range("AE6").Value = 99465
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 800
$ws.Range("AH6").Value = 1.52
$ws.Range("AI6").Value = 30.93
$ws.Range("AJ6").Value = 104909237

# Row 7
$ws.Range("D7").Value = 90072
$ws.Range("E7").Value = 1950
$ws.Range("G7").Value = 5880
$ws.Range("H7").Value = 6690
$ws.Range("I7").Value = 7202
$ws.Range("K7").Value = 296440
$ws.Range("L7").Value = 206270
$ws.Range("M7").Value = 90170
$ws.Range("N7").Value = 76600
$ws.Range("P7").Value = 210
$ws.Range("R7").Value = -9620
$ws.Range("S7").Value = 11660
$ws.Range("T7").Value = 4290
$ws.Range("W7").Value = 2.17
$ws.Range("X7").Value = 7.43
$ws.Range("Y7").Value = 9.75
$ws.Range("Z7").Value = 2.34
$ws.Range("AA7").Value = 228.76
$ws.Range("AC7").Value = 6801
$ws.Range("AD7").Value = 5.19
$ws.Range("AE7").Value = 107114
$ws.Range("AF7").Value = 0.33
$ws.Range("AG7").Value = 1200
$ws.Range("AH7").Value = 3.4
$ws.Range("AI7").Value = 17.48
$ws.Range("Q7").ClearContents()  # cell removed entirely in target
$ws.Range("U7").ClearContents()  # cell removed entirely in target

# Row 8
$ws.Range("D8").Value = 93092
$ws.Range("E8").Value = 2064
$ws.Range("G8").Value = 5170
$ws.Range("H8").Value = 4290
$ws.Range("I8").Value = 3322
$ws.Range("K8").Value = 308200
$ws.Range("L8").Value = 214580
$ws.Range("M8").Value = 93620
$ws.Range("N8").Value = 79810
$ws.Range("P8").Value = 210
$ws.Range("Q8").Value = 11900
$ws.Range("R8").Value = -6360
$ws.Range("S8").Value = 1640
$ws.Range("T8").Value = 4250
$ws.Range("U8").Value = 7650
$ws.Range("W8").Value = 2.22
$ws.Range("X8").Value = 4.61
$ws.Range("Y8").Value = 4.25
$ws.Range("Z8").Value = 1.42
$ws.Range("AA8").Value = 229.2
$ws.Range("AC8").Value = 3137
$ws.Range("AD8").Value = 11.25
$ws.Range("AE8").Value = 111602
$ws.Range("AF8").Value = 0.32
$ws.Range("AG8").Value = 1200
$ws.Range("AH8").Value = 3.4
$ws.Range("AI8").Value = 37.9

# Row 9
$ws.Range("D9").Value = 96490
$ws.Range("E9").Value = 2280
$ws.Range("G9").Value = 3830
$ws.Range("H9").Value = 3180
$ws.Range("I9").Value = 2890
$ws.Range("K9").Value = 318150
$ws.Range("L9").Value = 222200
$ws.Range("M9").Value = 95950
$ws.Range("N9").Value = 81850
$ws.Range("P9").Value = 210
$ws.Range("Q9").Value = 10790
$ws.Range("R9").Value = -6560
$ws.Range("S9").Value = 530
$ws.Range("T9").Value = 4250
$ws.Range("U9").Value = 6540
$ws.Range("W9").Value = 2.36
$ws.Range("X9").Value = 3.3
$ws.Range("Y9").Value = 3.57
$ws.Range("Z9").Value = 1.01
$ws.Range("AA9").Value = 231.58
$ws.Range("AC9").Value = 2729
$ws.Range("AD9").Value = 12.93
$ws.Range("AE9").Value = 114455
$ws.Range("AF9").Value = 0.31
$ws.Range("AG9").Value = 1200
$ws.Range("AH9").Value = 3.4
$ws.Range("AI9").Value = 43.56
